$wb = $excel.ActiveWorkbook

# The data lives on the second sheet ("dc.contributor.author")
$ws = $wb.Worksheets.Item("dc.contributor.author")

# Replace the old "::" separator notation with the new "$$" separator
# notation used to reference objects by their business identifiers.
$ws.Range("B2").Value = "Author1`$`$authority1"
$ws.Range("C3").Value = "OrgUnit2`$`$authority2`$`$400"

# Move the active selection to B3 (this sheet becomes the active sheet,
# matching the selection recorded in the saved workbook).
$ws.Activate()
$ws.Range("B3").Select()
